$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so formatted numeric strings
# (e.g. trailing zeros, thousands-dot formatting) survive the write.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '70.081.31'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.477.76'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '613.40'
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").Value = '168.54'
$ws.Range("E6").Value = '  -3.07%  '
$ws.Range("D7").Value = '3.474.26'
$ws.Range("E7").Value = '  -2.11%  '
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  -2.55%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '0.194'
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -2.60%  '
$ws.Range("D12").Value = '0.567'
$ws.Range("E12").Value = '  -3.32%  '
$ws.Range("D13").Value = '44.65'
$ws.Range("E13").Value = '  -3.81%  '
$ws.Range("D14").Value = '0.0000269'
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").Value = '4.032.61'
$ws.Range("E15").Value = '  -2.29%  '
$ws.Range("D16").Value = '8.23'
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '590.62'
$ws.Range("E17").Value = '  -2.87%  '
$ws.Range("D18").Value = '3.493.14'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("D19").Value = '70.098.98'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").Value = '17.25'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("D22").Value = '0.858'
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("D23").Value = '8.80'
$ws.Range("E23").Value = '  -5.24%  '
$ws.Range("D24").Value = '96.03'
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = '15.23'
$ws.Range("E25").Value = '  -3.13%  '
$ws.Range("D26").Value = '3.62'
$ws.Range("E26").Value = '  -3.00%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '2.49'
$ws.Range("E28").Value = '  -4.47%  '
$ws.Range("D29").Value = '33.13'
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("D30").Value = '8.69'
$ws.Range("E30").Value = '  -4.14%  '
$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -3.98%  '
$ws.Range("D32").Value = '2.84'
$ws.Range("E32").Value = '  -7.02%  '
$ws.Range("E33").Value = '  -2.88%  '
$ws.Range("D34").Value = '6.60'
$ws.Range("E34").Value = '  -6.12%  '
$ws.Range("D35").Value = '574.66'
$ws.Range("E35").Value = '  -21.33%  '
$ws.Range("D36").Value = '10.68'
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("D37").Value = '0.0483'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '0.0965'
$ws.Range("E38").Value = '  -4.01%  '
$ws.Range("E39").Value = '  +0.27%  '
$ws.Range("D40").Value = '56.35'
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("E41").Value = '  -0.96%  '
$ws.Range("D42").Value = '3.22'
$ws.Range("E42").Value = '  -9.97%  '
$ws.Range("D43").Value = '3.269.97'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("D44").Value = '0.0₃0700'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  -5.77%  '
$ws.Range("D46").Value = '30.97'
$ws.Range("E46").Value = '  -4.89%  '
$ws.Range("E47").Value = '  -5.67%  '
$ws.Range("D48").Value = '2.41'
$ws.Range("E48").Value = '  -6.79%  '
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").Value = '133.61'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("E51").Value = '  -0.01%  '
